$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value2 = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "dnasr281@gmail.com, admin@admin.com") {
        $cell.Value2 = "admin@admin.com, dnasr281@gmail.com"
    }
    elseif ($val -eq "system, System, backup@backdoor.com") {
        $cell.Value2 = "System, system, backup@backdoor.com"
    }
}
